$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (learning rate 0.0001 run -> became first entry, 0.005 swapped to I2) ---
$ws.Range("B2").Value = 0.0001
$ws.Range("D2").Value = 1000
$ws.Range("I2").Value = 0.005
$ws.Range("K2").Value = 0.9932000041007996
$ws.Range("L2").Value = 0.9843000173568726
$ws.Range("M2").Value = 123.85
$ws.Range("N2").Value = 0.008
$ws.Range("O2").Value = 0.007900000000000001
$ws.Range("P2").Value = 16
$ws.Range("Q2").Value = 7.7406
$ws.Range("R2").Value = 0.9922999739646912
$ws.Range("S2").Value = 0.9909999966621399
$ws.Range("T2").Value = 0.9909999966621399

# --- New row 3 ---
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 0.0001
$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 1000
$ws.Range("E3").Value = 200
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = "('tanh', 'relu')"
$ws.Range("H3").Value = 100
$ws.Range("I3").Value = 0.005
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 0.9937000274658203
$ws.Range("L3").Value = 0.9853000044822693
$ws.Range("M3").Value = 100.529
$ws.Range("N3").Value = 0.009900000000000001
$ws.Range("O3").Value = 0.0098
$ws.Range("P3").Value = 13
$ws.Range("Q3").Value = 7.733
$ws.Range("R3").Value = 0.9918000102043152
$ws.Range("S3").Value = 0.9915000200271606
$ws.Range("T3").Value = 0.9915000200271606

# --- New row 4 ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = 0.0001
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1000
$ws.Range("E4").Value = 200
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = "('tanh', 'relu')"
$ws.Range("H4").Value = 100
$ws.Range("I4").Value = 0.005
$ws.Range("J4").Value = 3
$ws.Range("K4").Value = 0.9926999807357788
$ws.Range("L4").Value = 0.9824000000953674
$ws.Range("M4").Value = 87.73999999999999
$ws.Range("N4").Value = 0.0113
$ws.Range("O4").Value = 0.0112
$ws.Range("P4").Value = 12
$ws.Range("Q4").Value = 7.3117
$ws.Range("R4").Value = 0.9911999702453613
$ws.Range("S4").Value = 0.9896000027656555
$ws.Range("T4").Value = 0.9896000027656555
